$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update claim number in F2 (NroSiniestro column) - trailing space preserved.
# Leading apostrophe forces text storage (not numeric) so the leading zero
# and trailing space survive, matching the cell's existing quote-prefix style.
$ws.Range("F2").Value = "'0420194406701 "

# Update the active cell selection to H7
$ws.Range("H7").Select()
